$wb = $excel.ActiveWorkbook

# Sheet "Overview": update Status cells for zh-cn (E2) and de-de (F2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E:F").ColumnWidth = 13.4101845877511

# Sheet "zh-cn": update Status cell C2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C:C").ColumnWidth = 13.4101845877511

# Sheet "de-de": update Status cell C2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C:C").ColumnWidth = 13.4101845877511
